# Y5_B2526_General_&_Special_Surgery_2_B1_reference_data.xlsx
# ---------------------------------------------------------------
# Re-create the workbook-level touch-up that Excel performs when the
# file is opened, the page setup is reset to the application defaults,
# and workbook protection metadata is (re)written on save.
#
# 1) Protect the workbook structure (adds the <workbookProtection/>
#    bookkeeping element that Excel writes once protection has been
#    touched).
# 2) Reset the worksheet's page margins to Excel's stock "Normal"
#    values (0.75"/0.75"/1"/1" for left/right/top/bottom and
#    0.5"/0.5" for header/footer) -- PageSetup margins are expressed
#    in points, so convert inches -> points (72 pt/in).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect the workbook structure (sheet order / add / delete / rename)
$wb.Protect($null, $true, $false)

# Reset page margins to Excel's default "Normal" page setup.
$ws.PageSetup.LeftMargin   = 0.75 * 72   # 54 pt  -> 0.75 in
$ws.PageSetup.RightMargin  = 0.75 * 72   # 54 pt  -> 0.75 in
$ws.PageSetup.TopMargin    = 1    * 72   # 72 pt  -> 1 in
$ws.PageSetup.BottomMargin = 1    * 72   # 72 pt  -> 1 in
$ws.PageSetup.HeaderMargin = 0.5  * 72   # 36 pt  -> 0.5 in
$ws.PageSetup.FooterMargin = 0.5  * 72   # 36 pt  -> 0.5 in
